# LOM3058.docx restructuring:
#   The "Objetivos" / "Programa resumido" / "Programa" / "Avaliação" /
#   "Bibliografia" sections had their paragraph bodies shuffled around
#   (and the evaluation criteria + bibliography text were rewritten/
#   relocated) while every paragraph's position, style and run
#   formatting (bold labels, italic EN text, ListBullet style, etc.)
#   stayed put. So: capture each distinct text block from its current
#   paragraph/run slot first, then write the blocks into their new
#   slots in a single pass - this sidesteps the cyclic reshuffle.

$d = $word.ActiveDocument

# --- capture original values from their current slots -------------
$OBJ_PT_SHORT = $d.Paragraphs(6).Range.Text.TrimEnd([char]13)
$OBJ_EN       = $d.Paragraphs(7).Range.Text.TrimEnd([char]13)
$DOCENTE      = $d.Paragraphs(9).Range.Text.TrimEnd([char]13)
$PROG_RES_PT  = $d.Paragraphs(11).Range.Text.TrimEnd([char]13)
$PROG_RES_EN  = $d.Paragraphs(12).Range.Text.TrimEnd([char]13)
$PROG_PT_LONG = $d.Paragraphs(14).Range.Text.TrimEnd([char]13)

$p17 = $d.Paragraphs(17)

$fMetodoVal = $p17.Range.Duplicate
$fMetodoVal.Find.Execute("Atividades avaliativas envolvendo o conteúdo teórico ministrado em sala de aula.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$METODO_VAL = $fMetodoVal.Text

$fCriterioVal = $p17.Range.Duplicate
$fCriterioVal.Find.Execute("Duas avaliações, sendo a primeira compondo 40% da nota e a segunda 60%. Os alunos que apresentarem média igual ou superior a 5 estarão aprovados, enquanto aqueles que tiverem média inferior a 3 estarão reprovados. Alunos com notas situadas no intervalo de 3 a 4,9 serão encaminhados à recuperação.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$CRITERIO_VAL = $fCriterioVal.Text

$fNormaVal = $p17.Range.Duplicate
$fNormaVal.Find.Execute("O aluno será submetido a um programa de estudos destinado a rever o conteúdo abordado na disciplina. Ao final deste período será aplicada uma nova avaliação. A nota final do aluno será a média aritmética desta avaliação com a nota anteriormente obtida, estando aprovados os alunos que tiverem nota final igual ou superior a 5.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$NORMA_VAL = $fNormaVal.Text

# Bibliography list (paragraph 19), 4 entries separated by line breaks
$BIB1 = 'G. ODIAN Principles of Polymerization, 3rd Edition, New York: Wiley-Interscience, 1991.'
$BIB2 = 'F. W. Billmeyer. Textbook of Polymer Chemistry, 3rd edition, New York: Wiley-Interscience, 1984.'
$BIB3 = 'C. E. Carraher. Introduction to Polymer Chemistry, 1st Edition, Boca Raton: Taylor and Francis, 2010.'
$BIB4 = 'S. V. Canevarolo. Ciência dos Polímeros: um texto básico para Engenheiros e Tecnólogos, 2ª. edição, São Paulo: Artliber, 2006.'
$LB = [char]11
$BIB_MERGED = "$BIB1$LB$BIB2$LB$BIB3$LB$BIB4"

# --- write captured values into their new slots --------------------
# (the 3 value-runs inside paragraph 17 form a rotation - Método's new
# value equals old Critério's value, Critério's new value equals old
# Norma's value, etc. - so writing the final text directly could make
# a later Find() match a just-written duplicate instead of the
# original slot. Stage through unique sentinel tokens first to avoid
# any such collision, then resolve the sentinels to real text.)

# Objetivos (PT) <- former "Programa resumido" PT text
$d.Paragraphs(6).Range.Text = $PROG_RES_PT
# Objetivos (EN, italic) <- former "Programa resumido" EN text
$d.Paragraphs(7).Range.Text = $PROG_RES_EN

# Docente(s) Responsável(eis) bullet <- former Objetivos PT text
$d.Paragraphs(9).Range.Text = $OBJ_PT_SHORT

# Programa resumido (PT) <- former "Programa" long PT text
$d.Paragraphs(11).Range.Text = $PROG_PT_LONG
# Programa resumido (EN, italic) <- former Objetivos EN text
$d.Paragraphs(12).Range.Text = $OBJ_EN

# Programa (PT, long) <- former "Método" value
$d.Paragraphs(14).Range.Text = $METODO_VAL
# Paragraph 15 (Programa EN long, italic) is unchanged.

# Stage 1: replace each of the 3 paragraph-17 value runs with a
# unique sentinel so none of the final strings can collide with one
# another while we still need to locate the *other* originals.
$TOK_METODO   = "@@TOKEN_METODO_VAL@@"
$TOK_CRITERIO = "@@TOKEN_CRITERIO_VAL@@"
$TOK_NORMA    = "@@TOKEN_NORMA_VAL@@"

$r = $p17.Range.Duplicate
$r.Find.Execute($METODO_VAL, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = $TOK_METODO

$r = $p17.Range.Duplicate
$r.Find.Execute($CRITERIO_VAL, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = $TOK_CRITERIO

$r = $p17.Range.Duplicate
$r.Find.Execute($NORMA_VAL, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = $TOK_NORMA

# Stage 2: resolve each sentinel to its real final text.
$r = $p17.Range.Duplicate
$r.Find.Execute($TOK_METODO, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = $CRITERIO_VAL

$r = $p17.Range.Duplicate
$r.Find.Execute($TOK_CRITERIO, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = $NORMA_VAL

$r = $p17.Range.Duplicate
$r.Find.Execute($TOK_NORMA, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = $BIB_MERGED

# Bibliografia list (paragraph 19) <- former Docente text
$d.Paragraphs(19).Range.Text = $DOCENTE
